$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: update four odd values (FlashScore re-scrape) ---
$ws.Range("G19").Value = 2.63
$ws.Range("AG19").Value = 7
$ws.Range("AJ19").Value = 29
$ws.Range("AN19").Value = 4.5

# --- Insert a brand-new fixture as row 39 (Barry vs Briton Ferry). ---
# This shifts the former rows 39-41 (Caernarfon/Newtown, Flint/Bala,
# Haverfordwest/Connahs Q.) down to rows 40-42 untouched, matching the
# dimension growing from A1:BD41 to A1:BD42.
$ws.Rows.Item(39).Insert()

$ws.Range("A39").Value = "x2szDPi5"
$ws.Range("B39").Value = "'12/10/2024"
$ws.Range("C39").Value = "10:30"
$ws.Range("D39").Value = "WALES - CYMRU PREMIER"
$ws.Range("E39").Value = "Barry"
$ws.Range("F39").Value = "Briton Ferry"
$ws.Range("G39").Value = 1.39
$ws.Range("H39").Value = 4.5
$ws.Range("I39").Value = 6.3
$ws.Range("J39").Value = 1.85
$ws.Range("K39").Value = 2.5
$ws.Range("L39").Value = 6
$ws.Range("M39").Value = 1.03
$ws.Range("N39").Value = 9
$ws.Range("O39").Value = 1.18
$ws.Range("P39").Value = 4.3
$ws.Range("Q39").Value = 1.55
$ws.Range("R39").Value = 2.3
$ws.Range("S39").Value = 1.29
$ws.Range("T39").Value = 3.3
$ws.Range("U39").Value = 1.75
$ws.Range("V39").Value = 1.95
$ws.Range("W39").Value = 8.5
$ws.Range("X39").Value = 7.5
$ws.Range("Y39").Value = 8.25
$ws.Range("Z39").Value = 9.75
$ws.Range("AA39").Value = 10.75
$ws.Range("AB39").Value = 22
$ws.Range("AC39").Value = 9
$ws.Range("AD39").Value = 9.25
$ws.Range("AE39").Value = 17.5
$ws.Range("AF39").Value = 70
$ws.Range("AG39").Value = 20
$ws.Range("AH39").Value = 45
$ws.Range("AI39").Value = 20
$ws.Range("AJ39").Value = 150
$ws.Range("AK39").Value = 65
$ws.Range("AL39").Value = 55
$ws.Range("AM39").Value = 450
$ws.Range("AN39").Value = 3.4
$ws.Range("AO39").Value = 6.3
$ws.Range("AP39").Value = 14.5
$ws.Range("AQ39").Value = 16.5
$ws.Range("AR39").Value = 40
$ws.Range("AS39").Value = 175
$ws.Range("AT39").Value = 3.3
$ws.Range("AU39").Value = 7.8
$ws.Range("AV39").Value = 65
$ws.Range("AW39").Value = 8
$ws.Range("AX39").Value = 35
$ws.Range("AY39").Value = 35
$ws.Range("AZ39").Value = 250
$ws.Range("BA39").Value = 250
$ws.Range("BB39").Value = 400
$ws.Range("BC39").Value = 51
$ws.Range("BD39").Value = 51


Write-Output "edit complete"
